$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 2236
$ws.Cells.Item(3, 6).Value = 265
$ws.Cells.Item(4, 6).Value = 162
$ws.Cells.Item(5, 6).Value = 157
$ws.Cells.Item(6, 6).Value = 297
$ws.Cells.Item(8, 6).Value = 669
$ws.Cells.Item(9, 6).Value = 493
$ws.Cells.Item(10, 6).Value = 603
$ws.Cells.Item(11, 6).Value = 363
$ws.Cells.Item(12, 6).Value = 58
$ws.Cells.Item(13, 6).Value = 338
$ws.Cells.Item(14, 6).Value = 944
$ws.Cells.Item(15, 6).Value = 205
$ws.Cells.Item(16, 6).Value = 125
$ws.Cells.Item(17, 6).Value = 83
$ws.Cells.Item(20, 6).Value = 204
$ws.Cells.Item(21, 6).Value = 80

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(8, 6).Value = 2422
$ws.Cells.Item(10, 6).Value = 14
$ws.Cells.Item(14, 6).Value = 94
$ws.Cells.Item(16, 6).Value = 2255

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(3, 6).Value = 325
$ws.Cells.Item(4, 6).Value = 157

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 2236
$ws.Cells.Item(6, 6).Value = 325
$ws.Cells.Item(7, 6).Value = 265
$ws.Cells.Item(8, 6).Value = 162
$ws.Cells.Item(9, 6).Value = 157
$ws.Cells.Item(10, 6).Value = 297
$ws.Cells.Item(15, 6).Value = 157
$ws.Cells.Item(16, 6).Value = 669
$ws.Cells.Item(17, 6).Value = 493
$ws.Cells.Item(18, 6).Value = 603
$ws.Cells.Item(19, 6).Value = 363
$ws.Cells.Item(20, 6).Value = 58
$ws.Cells.Item(21, 6).Value = 338
$ws.Cells.Item(22, 6).Value = 944
$ws.Cells.Item(24, 6).Value = 2423
$ws.Cells.Item(26, 6).Value = 14
$ws.Cells.Item(30, 6).Value = 205
$ws.Cells.Item(31, 6).Value = 125
$ws.Cells.Item(32, 6).Value = 83
$ws.Cells.Item(34, 6).Value = 94
$ws.Cells.Item(37, 6).Value = 204
$ws.Cells.Item(38, 6).Value = 80
$ws.Cells.Item(39, 6).Value = 2255
